$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 38.33049633333334
$ws.Range("H2").Value = 114.991489
$ws.Range("I2").Value = 0.5317874798120843
$ws.Range("J2").Value = 0.5317874798120843
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3676493333333333
$ws.Range("N2").Value = 1.102948
$ws.Range("O2").Value = 0.02558020875176611
$ws.Range("P2").Value = 0.02558020875176611
$ws.Range("Q2").Value = 14.09218142328578
$ws.Range("R2").Value = 126.829632809572
$ws.Range("S2").Value = 0.01360323474516872
$ws.Range("T2").Value = 0.01360323474516872

$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 38.33049633333334
$ws.Range("H3").Value = 114.991489
$ws.Range("I3").Value = 0.5317874798120843
$ws.Range("J3").Value = 0.5317874798120843
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.646212666666667
$ws.Range("N3").Value = 10.938638
$ws.Range("O3").Value = 0.2536952272455287
$ws.Range("P3").Value = 0.2536952272455287
$ws.Range("Q3").Value = 139.7611412502202
$ws.Range("R3").Value = 1257.850271251982
$ws.Range("S3").Value = 0.1349119455372537
$ws.Range("T3").Value = 0.1349119455372537

$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 38.33049633333334
$ws.Range("H4").Value = 114.991489
$ws.Range("I4").Value = 0.5317874798120843
$ws.Range("J4").Value = 0.5317874798120843
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.032567
$ws.Range("N4").Value = 0.097701
$ws.Range("O4").Value = 0.002265938172294887
$ws.Range("P4").Value = 0.002265938172294887
$ws.Range("Q4").Value = 1.248309274087667
$ws.Range("R4").Value = 11.234783466789
$ws.Range("S4").Value = 0.001204997550054698
$ws.Range("T4").Value = 0.001204997550054698

$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 38.33049633333334
$ws.Range("H5").Value = 114.991489
$ws.Range("I5").Value = 0.5317874798120843
$ws.Range("J5").Value = 0.5317874798120843
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.32598433333333
$ws.Range("N5").Value = 30.977953
$ws.Range("O5").Value = 0.7184586258304102
$ws.Range("P5").Value = 0.7184586258304102
$ws.Range("Q5").Value = 395.8001046268908
$ws.Range("R5").Value = 3562.200941642017
$ws.Range("S5").Value = 0.3820673019796071
$ws.Range("T5").Value = 0.3820673019796071

$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 14.943524
$ws.Range("H6").Value = 44.830572
$ws.Range("I6").Value = 0.2073226210890634
$ws.Range("J6").Value = 0.2073226210890634
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3676493333333333
$ws.Range("N6").Value = 1.102948
$ws.Range("O6").Value = 0.02558020875176611
$ws.Range("P6").Value = 0.02558020875176611
$ws.Range("Q6").Value = 5.493976636250667
$ws.Range("R6").Value = 49.445789726256
$ws.Range("S6").Value = 0.005303355926421548
$ws.Range("T6").Value = 0.005303355926421548

$ws.Range("D7").Value = "FAPs"
$ws.Range("G7").Value = 14.943524
$ws.Range("H7").Value = 44.830572
$ws.Range("I7").Value = 0.2073226210890634
$ws.Range("J7").Value = 0.2073226210890634
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.646212666666667
$ws.Range("N7").Value = 10.938638
$ws.Range("O7").Value = 0.2536952272455287
$ws.Range("P7").Value = 0.2536952272455287
$ws.Range("Q7").Value = 54.48726649343735
$ws.Range("R7").Value = 490.3853984409361
$ws.Range("S7").Value = 0.05259675947032857
$ws.Range("T7").Value = 0.05259675947032857

$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("G8").Value = 14.943524
$ws.Range("H8").Value = 44.830572
$ws.Range("I8").Value = 0.2073226210890634
$ws.Range("J8").Value = 0.2073226210890634
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.032567
$ws.Range("N8").Value = 0.097701
$ws.Range("O8").Value = 0.002265938172294887
$ws.Range("P8").Value = 0.002265938172294887
$ws.Range("Q8").Value = 0.4866657461080001
$ws.Range("R8").Value = 4.379991714972
$ws.Range("S8").Value = 0.0004697802411059376
$ws.Range("T8").Value = 0.0004697802411059376

$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 14.943524
$ws.Range("H9").Value = 44.830572
$ws.Range("I9").Value = 0.2073226210890634
$ws.Range("J9").Value = 0.2073226210890634
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.32598433333333
$ws.Range("N9").Value = 30.977953
$ws.Range("O9").Value = 0.7184586258304102
$ws.Range("P9").Value = 0.7184586258304102
$ws.Range("Q9").Value = 154.3065947087907
$ws.Range("R9").Value = 1388.759352379116
$ws.Range("S9").Value = 0.1489527254512073
$ws.Range("T9").Value = 0.1489527254512073

$ws.Range("D10").Value = "ECs"
$ws.Range("G10").Value = 8.167063666666666
$ws.Range("H10").Value = 24.501191
$ws.Range("I10").Value = 0.1133077476219524
$ws.Range("J10").Value = 0.1133077476219524
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3676493333333333
$ws.Range("N10").Value = 1.102948
$ws.Range("O10").Value = 0.02558020875176611
$ws.Range("P10").Value = 0.02558020875176611
$ws.Range("Q10").Value = 3.002615512340888
$ws.Range("R10").Value = 27.023539611068
$ws.Range("S10").Value = 0.002898435837361974
$ws.Range("T10").Value = 0.002898435837361974

$ws.Range("D11").Value = "FAPs"
$ws.Range("G11").Value = 8.167063666666666
$ws.Range("H11").Value = 24.501191
$ws.Range("I11").Value = 0.1133077476219524
$ws.Range("J11").Value = 0.1133077476219524
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.646212666666667
$ws.Range("N11").Value = 10.938638
$ws.Range("O11").Value = 0.2536952272455287
$ws.Range("P11").Value = 0.2536952272455287
$ws.Range("Q11").Value = 29.77885099087311
$ws.Range("R11").Value = 268.009658917858
$ws.Range("S11").Value = 0.02874563478163024
$ws.Range("T11").Value = 0.02874563478163025

$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 8.167063666666666
$ws.Range("H12").Value = 24.501191
$ws.Range("I12").Value = 0.1133077476219524
$ws.Range("J12").Value = 0.1133077476219524
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.032567
$ws.Range("N12").Value = 0.097701
$ws.Range("O12").Value = 0.002265938172294887
$ws.Range("P12").Value = 0.002265938172294887
$ws.Range("Q12").Value = 0.2659767624323333
$ws.Range("R12").Value = 2.393790861891
$ws.Range("S12").Value = 0.0002567483505533372
$ws.Range("T12").Value = 0.0002567483505533373

$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 8.167063666666666
$ws.Range("H13").Value = 24.501191
$ws.Range("I13").Value = 0.1133077476219524
$ws.Range("J13").Value = 0.1133077476219524
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.32598433333333
$ws.Range("N13").Value = 30.977953
$ws.Range("O13").Value = 0.7184586258304102
$ws.Range("P13").Value = 0.7184586258304102
$ws.Range("Q13").Value = 84.33297147133587
$ws.Range("R13").Value = 758.996743242023
$ws.Range("S13").Value = 0.08140692865240688
$ws.Range("T13").Value = 0.08140692865240688

$ws.Range("D14").Value = "ECs"
$ws.Range("G14").Value = 5.834252333333334
$ws.Range("H14").Value = 17.502757
$ws.Range("I14").Value = 0.08094292121735479
$ws.Range("J14").Value = 0.08094292121735479
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.3676493333333333
$ws.Range("N14").Value = 1.102948
$ws.Range("O14").Value = 0.02558020875176611
$ws.Range("P14").Value = 0.02558020875176611
$ws.Range("Q14").Value = 2.144958980848445
$ws.Range("R14").Value = 19.304630827636
$ws.Range("S14").Value = 0.002070536821717694
$ws.Range("T14").Value = 0.002070536821717694

$ws.Range("D15").Value = "FAPs"
$ws.Range("G15").Value = 5.834252333333334
$ws.Range("H15").Value = 17.502757
$ws.Range("I15").Value = 0.08094292121735479
$ws.Range("J15").Value = 0.08094292121735479
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.646212666666667
$ws.Range("N15").Value = 10.938638
$ws.Range("O15").Value = 0.2536952272455287
$ws.Range("P15").Value = 0.2536952272455287
$ws.Range("Q15").Value = 21.27292475832956
$ws.Range("R15").Value = 191.456322824966
$ws.Range("S15").Value = 0.02053483279215375
$ws.Range("T15").Value = 0.02053483279215375

$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 5.834252333333334
$ws.Range("H16").Value = 17.502757
$ws.Range("I16").Value = 0.08094292121735479
$ws.Range("J16").Value = 0.08094292121735479
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.032567
$ws.Range("N16").Value = 0.097701
$ws.Range("O16").Value = 0.002265938172294887
$ws.Range("P16").Value = 0.002265938172294887
$ws.Range("Q16").Value = 0.1900040957396667
$ws.Range("R16").Value = 1.710036861657
$ws.Range("S16").Value = 0.0001834116549634619
$ws.Range("T16").Value = 0.0001834116549634619

$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 5.834252333333334
$ws.Range("H17").Value = 17.502757
$ws.Range("I17").Value = 0.08094292121735479
$ws.Range("J17").Value = 0.08094292121735479
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 10.32598433333333
$ws.Range("N17").Value = 30.977953
$ws.Range("O17").Value = 0.7184586258304102
$ws.Range("P17").Value = 0.7184586258304102
$ws.Range("Q17").Value = 60.24439819071345
$ws.Range("R17").Value = 542.1995837164211
$ws.Range("S17").Value = 0.05815413994851987
$ws.Range("T17").Value = 0.05815413994851987

$ws.Range("D18").Value = "ECs"
$ws.Range("G18").Value = 4.803262333333334
$ws.Range("H18").Value = 14.409787
$ws.Range("I18").Value = 0.06663923025954499
$ws.Range("J18").Value = 0.066639230259545
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.3676493333333333
$ws.Range("N18").Value = 1.102948
$ws.Range("O18").Value = 0.02558020875176611
$ws.Range("P18").Value = 0.02558020875176611
$ws.Range("Q18").Value = 1.765916194675111
$ws.Range("R18").Value = 15.893245752076
$ws.Range("S18").Value = 0.00170464542109617
$ws.Range("T18").Value = 0.00170464542109617

$ws.Range("D19").Value = "FAPs"
$ws.Range("G19").Value = 4.803262333333334
$ws.Range("H19").Value = 14.409787
$ws.Range("I19").Value = 0.06663923025954499
$ws.Range("J19").Value = 0.066639230259545
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 3.646212666666667
$ws.Range("N19").Value = 10.938638
$ws.Range("O19").Value = 0.2536952272455287
$ws.Range("P19").Value = 0.2536952272455287
$ws.Range("Q19").Value = 17.51371596112289
$ws.Range("R19").Value = 157.623443650106
$ws.Range("S19").Value = 0.01690605466416238
$ws.Range("T19").Value = 0.01690605466416239

$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 4.803262333333334
$ws.Range("H20").Value = 14.409787
$ws.Range("I20").Value = 0.06663923025954499
$ws.Range("J20").Value = 0.066639230259545
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.032567
$ws.Range("N20").Value = 0.097701
$ws.Range("O20").Value = 0.002265938172294887
$ws.Range("P20").Value = 0.002265938172294887
$ws.Range("Q20").Value = 0.1564278444096667
$ws.Range("R20").Value = 1.407850599687
$ws.Range("S20").Value = 0.0001510003756174515
$ws.Range("T20").Value = 0.0001510003756174515

$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 4.803262333333334
$ws.Range("H21").Value = 14.409787
$ws.Range("I21").Value = 0.06663923025954499
$ws.Range("J21").Value = 0.066639230259545
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 10.32598433333333
$ws.Range("N21").Value = 30.977953
$ws.Range("O21").Value = 0.7184586258304102
$ws.Range("P21").Value = 0.7184586258304102
$ws.Range("Q21").Value = 49.59841160289011
$ws.Range("R21").Value = 446.385704426011
$ws.Range("S21").Value = 0.04787752979866898
$ws.Range("T21").Value = 0.04787752979866899
